# Update 2p0. Convention change to support multi-axle vehicles
#
# Renames the shared labels "sAxleF" -> "sAxle1" and "sAxleR" -> "sAxle2"
# (cells A5 and A6 respectively) on every vehicle-data worksheet, and
# restores the view state so that the "Sedan_HambaLG" sheet is the active /
# selected tab (matching the author's final save state).

$wb = $excel.ActiveWorkbook

# --- 1. Rename the axle labels on every worksheet -------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Range("A5").Value = "sAxle1"
    $ws.Range("A6").Value = "sAxle2"
}

# --- 2. Restore per-sheet selection / view state ---------------------------

# Sheet "Sedan_Hamba": selection on the axle label rows, scrolled to A17
$wsHamba = $wb.Worksheets.Item("Sedan_Hamba")
$wsHamba.Activate()
$wsHamba.Range("A5:A6").Select()
$wsHamba.Range("A17").Select()

# Sheet "Bus_Makhulu": selection on the axle label rows, scrolled to A2
$wsBus = $wb.Worksheets.Item("Bus_Makhulu")
$wsBus.Activate()
$wsBus.Range("A5:A6").Select()
$wsBus.Range("A2").Select()

# Sheet "Sedan_HambaLG": becomes the active/selected tab, scrolled to A2
$wsHambaLG = $wb.Worksheets.Item("Sedan_HambaLG")
$wsHambaLG.Activate()
$wsHambaLG.Range("A2").Select()
